# "add words in June 7th"
#
# The paragraph that carries the (hidden) "_GoBack" bookmark -- the empty
# paragraph right after the June 7th diary entry -- gets a new run of text
# ("改的第一遍") inserted in front of the bookmark, and a brand-new empty
# paragraph is added immediately afterwards.
#
# We locate that paragraph, then replace it (and splice in the new empty
# paragraph after it) via Range.InsertXML using a minimal Flat-OPC package,
# so that the run/paragraph formatting (rFonts hint + lang) matches exactly
# what a real Word edit in that position would produce.

$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r") {
        $prev = $null
        if ($i -gt 1) { $prev = $d.Paragraphs.Item($i - 1) }
        if ($prev -ne $null -and $prev.Range.Text -like "*高考第一天*") {
            $target = $p
            break
        }
    }
}

if ($target -eq $null) {
    throw "could not locate the target paragraph (the one holding _GoBack)"
}

$r = $target.Range
$r.Collapse(1)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:hint="default"/>
                <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
              </w:rPr>
              <w:t>改的第一遍</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:hint="default"/>
                <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
